$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Cases")

# Delete rows that correspond to removed test cases (Search45, Search46,
# Search50-Search53 and their associated OPQA rows) - highest row numbers
# first so earlier row numbers remain valid.
$ws.Rows("51:54").Delete()
$ws.Rows("46:47").Delete()

# Update the active selection to match the post-edit workbook state.
$ws.Range("A49").Select()
